$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Record Days")

# The Date column in this sheet stores plain text values such as
# "2026-02-18" (not real Excel dates). Briefly mark the target cells as
# Text so the assigned strings aren't auto-converted into date serial
# numbers, then restore the default "Normal" style so no stray
# number-format style is left attached to the cells.
$ws.Range("B18:B19").NumberFormat = "@"

# New row 18 content (previously held by row 19): no RxNumber value.
$ws.Range("A18").Value = "628994d9-3a33-43b2-9214-c3bd9793ccd4"
$ws.Range("B18").Value = "2026-02-18"
$ws.Range("C18").Value = ""
$ws.Range("D18").Value = "draft"

# New row 19 content (previously held by row 18), with the RxNumber
# label corrected to match the sheet's naming convention.
$ws.Range("A19").Value = "e432f0fe-1383-44a2-990c-5f787da5008a"
$ws.Range("B19").Value = "2026-02-04"
$ws.Range("C19").Value = "RX EP 1 - 5"
$ws.Range("D19").Value = "ready"

$ws.Range("B18:B19").Style = "Normal"
